$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (student id 110): name "h" -> "Max", email -> max@x.com
$ws.Range("B11").Value = "Max"
$ws.Range("C11").Value = "max@x.com"
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:max@x.com")
$ws.Range("C11").Style = "Normal"
$ws.Range("C11").Style = "Hyperlink"

# Row 12 (student id 111): name "i" -> "Kat", email -> kat@x.com
$ws.Range("B12").Value = "Kat"
$ws.Range("C12").Value = "kat@x.com"
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:kat@x.com")
$ws.Range("C12").Style = "Normal"
$ws.Range("C12").Style = "Hyperlink"

$ws.Range("B13").Select()
